# Applies the "Tested stimulated stock. Fixed bug with calculating appreciation" edit.
#
# Summary of changes:
#  1. Data sheet: header K1 renamed invested_capital -> capital_invested; the
#     "month" column (A) is now 1-based instead of 0-based, and the
#     capital_gains_month / equity / asset_value / effective_income /
#     annual_ROI_extrapolated columns (F, I, J, L, M) now hold a compounding
#     market-growth calculation instead of the old flat placeholder values.
#  2. Revenue Ledger: every month block now carries a "Dividend" label/value
#     pair (previously only some months had one) - label text also renamed
#     from "Dividends" to "Dividend".
#  3. Investments Ledger: a new row 2 is added with a "Dividend Reinvestment"
#     label/value pair for every month.
#  4. Capital Gains Ledger: the label changes from "Capital Appreciation" to
#     "Market Growth", and the flat 1.0 placeholder values are replaced with
#     the same compounding market-growth values used on the Data sheet.

$wb = $excel.ActiveWorkbook

$dataSheet = $wb.Worksheets.Item("Data")
$revenueSheet = $wb.Worksheets.Item("Revenue Ledger")
$investSheet = $wb.Worksheets.Item("Investments Ledger")
$capGainsSheet = $wb.Worksheets.Item("Capital Gains Ledger")

# ---------------------------------------------------------------------------
# 1. Data sheet
# ---------------------------------------------------------------------------

# Header rename: invested_capital -> capital_invested (column K, row 1)
$dataSheet.Cells.Item(1, 11).Value = "capital_invested"

# Per-month (row 2..15) new values for columns A, F, I, J, L, M
$monthCol = @(1.0, 2.0, 3.0, 4.0, 5.0, 6.0, 7.0, 8.0, 9.0, 10.0, 11.0, 12.0, 13.0, 14.0)
$gainsCol = @(0.0797414043, 0.08101314261194721, 0.0823051629636738, 0.08361778881884652, 0.08495134879982165, 0.08630617676991702, 0.08768261191699647, 0.0890809988383871, 0.09050168762715104, 0.09194503395973284, 0.09341139918500489, 0.09490115041473293, 0.09641466061548419, 0.09795230870200153)
$equityCol = @(10.0797414043, 10.160754546911948, 10.243059709875622, 10.326677498694469, 10.41162884749429, 10.497935024264207, 10.585617636181203, 10.67469863501959, 10.765200322646741, 10.857145356606475, 10.95055675579148, 11.045457906206213, 11.141872566821696, 11.239824875523698)
$roiCol = @(0.09917416475042096, 0.09998692994627234, 0.10080004237548734, 0.10161340020425058, 0.10242690147587585, 0.10324044416197853, 0.10405392621369414, 0.10486724561293181, 0.10568030042360554, 0.10649298884274239, 0.10730520925149789, 0.10811686026592171, 0.1089278407875176, 0.10973805005346815)

for ($i = 0; $i -lt 14; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 1).Value = $monthCol[$i]   # A: month
    $dataSheet.Cells.Item($row, 6).Value = $gainsCol[$i]   # F: capital_gains_month
    $dataSheet.Cells.Item($row, 9).Value = $equityCol[$i]  # I: equity
    $dataSheet.Cells.Item($row, 10).Value = $equityCol[$i] # J: asset_value
    $dataSheet.Cells.Item($row, 12).Value = $gainsCol[$i]  # L: effective_income
    $dataSheet.Cells.Item($row, 13).Value = $roiCol[$i]    # M: annual_ROI_extrapolated
}

# ---------------------------------------------------------------------------
# Column layout shared by the ledger sheets: one (label, value) pair of
# columns per month, 14 months (Month: 0 .. Month: 13).
# ---------------------------------------------------------------------------
$labelCols = @(1, 4, 7, 10, 13, 16, 19, 22, 25, 28, 31, 34, 37, 40)
$valueCols = @(2, 5, 8, 11, 14, 17, 20, 23, 26, 29, 32, 35, 38, 41)

# ---------------------------------------------------------------------------
# 2. Revenue Ledger: fill in a "Dividend" label/value pair for every month.
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt 14; $i++) {
    $revenueSheet.Cells.Item(2, $labelCols[$i]).Value = "Dividend"
    $revenueSheet.Cells.Item(2, $valueCols[$i]).Value = 0.0
}

# ---------------------------------------------------------------------------
# 3. Investments Ledger: add row 2 with "Dividend Reinvestment" pairs.
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt 14; $i++) {
    $investSheet.Cells.Item(2, $labelCols[$i]).Value = "Dividend Reinvestment"
    $investSheet.Cells.Item(2, $valueCols[$i]).Value = 0.0
}

# ---------------------------------------------------------------------------
# 4. Capital Gains Ledger: rename label to "Market Growth" and set the real
#    compounding market-growth values (mirrors the Data sheet F column).
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt 14; $i++) {
    $capGainsSheet.Cells.Item(2, $labelCols[$i]).Value = "Market Growth"
    $capGainsSheet.Cells.Item(2, $valueCols[$i]).Value = $gainsCol[$i]
}
